$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TPM-derived values on row 2 to reflect the new TPM computation.
$ws.Range("G2").Value = 1.524170333333333
$ws.Range("H2").Value = 4.572511
$ws.Range("Q2").Value = 0.014764638019
$ws.Range("R2").Value = 0.132881742171
